# Replace double quotes with single quotes in the English (en_US) column
# of select dialogue lines, per commit "update on 20210731 画中人".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = "[name=`"Capone`"]   Listen up, once the so-called 'boss' opens fire, that's the official announcement of the end of your affiliation with the Gambino Family. Don't go easy on 'em.`n"

$ws.Range("C57").Value = "[name=`"Rat King`"]   Put down your weapons, kneel before me, and your 'Family' will have at least some chance to survive.`n"

$ws.Range("C66").Value = "[name=`"Rat King`"]   Look around you at your companions. This city has bled so much already. I really don't want to see 'Family' killing one another on the program.`n"

$ws.Range("C85").Value = "[name=`"Rat King`"]   'Your Family.'`n"

$ws.Range("C116").Value = "[name=`"Capone`"]   That's because I too am 'Sicilian.'`n"

$ws.Range("C117").Value = "[name=`"Capone`"]   Or I guess I should ask if you really thought you could wipe out one of Siracusa's 'Families' and get away with it?`n"

# Writing long strings into these cells nudges Excel's auto row-height
# calculation away from the sheet's default. Re-run AutoFit so the row
# metadata round-trips unchanged (matches upstream diff, which only
# touches the text content, not row heights).
foreach ($r in 11, 57, 66, 85, 116, 117) {
    $ws.Rows.Item($r).AutoFit() | Out-Null
}
